$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in row 26 (UMASS)
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 1

# Update the active cell selection to E26
$ws.Range("E26").Select()
